$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 265
$ws1.Range("F4").Value = 276
$ws1.Range("F5").Value = 2884
$ws1.Range("F8").Value = 2237
$ws1.Range("F9").Value = 2237
$ws1.Range("F10").Value = 1411
$ws1.Range("F12").Value = 445
$ws1.Range("F13").Value = 87
$ws1.Range("F14").Value = 2574
$ws1.Range("F16").Value = 1388
$ws1.Range("F17").Value = 4813
$ws1.Range("F19").Value = 5303
$ws1.Range("F20").Value = 5303
$ws1.Range("F21").Value = 1859
$ws1.Range("F22").Value = 2929
$ws1.Range("F23").Value = 3336
$ws1.Range("F24").Value = 188
$ws1.Range("F25").Value = 1598
$ws1.Range("F28").Value = 122
$ws1.Range("F30").Value = 315
$ws1.Range("F32").Value = 2048
$ws1.Range("F33").Value = 1
$ws1.Range("F34").Value = 124
$ws1.Range("F35").Value = 301
$ws1.Range("F36").Value = 762
$ws1.Range("F39").Value = 435

# Sheet "演出" (Performances) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 25
$ws2.Range("F18").Value = 51

# Sheet "全部类型" (All Types) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 265
$ws4.Range("F9").Value = 276
$ws4.Range("F10").Value = 2884
$ws4.Range("F12").Value = 2237
$ws4.Range("F13").Value = 2237
$ws4.Range("F14").Value = 1411
$ws4.Range("F17").Value = 445
$ws4.Range("F18").Value = 87
$ws4.Range("F19").Value = 25
$ws4.Range("F20").Value = 2574
$ws4.Range("F21").Value = 1388
$ws4.Range("F25").Value = 4813
$ws4.Range("F27").Value = 5303
$ws4.Range("F28").Value = 5303
$ws4.Range("F29").Value = 1859
$ws4.Range("F30").Value = 2929
$ws4.Range("F31").Value = 3336
$ws4.Range("F33").Value = 188
$ws4.Range("F36").Value = 1598
$ws4.Range("F40").Value = 122
$ws4.Range("F42").Value = 315
$ws4.Range("F43").Value = 51
$ws4.Range("F44").Value = 2048
$ws4.Range("F45").Value = 1
$ws4.Range("F46").Value = 124
$ws4.Range("F47").Value = 301
$ws4.Range("F48").Value = 762
$ws4.Range("F51").Value = 435
